# Fix sheet removed from spreadsheet
# The 4th sheet was generically named "Sheet1" instead of following the
# date-based naming convention used by the other sheets (110816, 110916,
# 111016). Rename it to "111116" to restore it properly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "111116"

$ws3 = $wb.Worksheets.Item("111016")
$ws3.Range("I6").Select()

$ws4 = $wb.Worksheets.Item("111116")
$ws4.Range("M20").Select()
